$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "Data utworzenia: 2025-08-05 21:28:11" "Data utworzenia: 2025-08-06 23:04:08"

Replace-Text " Odkryj tajemnicę obliczania środka ciężkości skomplikowanych figur geometrycznych! Dowiedz się, jak dzielić figurę na proste kształty, znajdować ich środki ciężkości i używać tych informacji do obliczenia środka ciężkości całej figury. Idealne dla studentów i profesjonalistów! " " Odkryj tajemnicę obliczania środka ciężkości dla złożonych figur płaskich! Krok po kroku przeprowadzimy Cię przez proces dzielenia figur na prostsze kształty i obliczania ich środków ciężkości. Idealne dla studentów i profesjonalistów! "

Replace-Text "[PIN] #Matematyka [PIN] #Geometria [PIN] #Edukacja [PIN] #Nauka" "[PIN] #środekciężkości [PIN] #matematyka [PIN] #nauka"

Replace-Text " Uncover the secret of calculating the center of gravity for complicated geometric figures! Learn how to divide a figure into simple shapes, find their centers of gravity and use this information to calculate the center of gravity of the entire figure. Perfect for students and professionals! " " Discover the secret of calculating the center of gravity for complex flat figures! Step by step, we will guide you through the process of dividing figures into simpler shapes and calculating their centers of gravity. Perfect for students and professionals! "

Replace-Text "[PIN] #Mathematics [PIN] #Geometry [PIN] #Education [PIN] #Science" "[PIN] #centerofgravity [PIN] #mathematics [PIN] #learning"

Replace-Text "środek ciężkości, figury geometryczne, edukacja, nauka" "center of gravity, complex figures"
